# Fix #9722 - [Feature] Translate export search reports
#
# Renames the "Fournisseurs"/"Usines" sheets to "Suppliers"/"Factories" and
# translates the remaining French column headers (on the Factories and
# Contacts sheets) to English. Also refreshes the header-row fill color.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename sheets (definedNames referencing the sheet names follow
#    automatically).
# ---------------------------------------------------------------------
$wsSuppliers = $wb.Worksheets.Item("Fournisseurs")
$wsSuppliers.Name = "Suppliers"

$wsFactories = $wb.Worksheets.Item("Usines")
$wsFactories.Name = "Factories"

$wsContacts = $wb.Worksheets.Item("Contacts")

# ---------------------------------------------------------------------
# 2. Translate the "Factories" header row (B3:M3) from French to English.
# ---------------------------------------------------------------------
$wsFactories.Range("B3").Value = "Name"
$wsFactories.Range("C3").Value = "Title"
$wsFactories.Range("D3").Value = "Phone"
$wsFactories.Range("E3").Value = "Address 1"
$wsFactories.Range("F3").Value = "Address 2"
$wsFactories.Range("G3").Value = "Address 3"
$wsFactories.Range("H3").Value = "Postal code"
$wsFactories.Range("I3").Value = "City"
$wsFactories.Range("J3").Value = "Country"
$wsFactories.Range("K3").Value = "Packager Code"
$wsFactories.Range("L3").Value = "Certifications"
$wsFactories.Range("M3").Value = "Agreement number"

# ---------------------------------------------------------------------
# 3. Translate the "Contacts" header row (B3:I3) from French to English.
# ---------------------------------------------------------------------
$wsContacts.Range("B3").Value = "First name"
$wsContacts.Range("C3").Value = "Last name"
$wsContacts.Range("D3").Value = "Email"
$wsContacts.Range("E3").Value = "Posting"
$wsContacts.Range("F3").Value = "Office phone"
$wsContacts.Range("G3").Value = "Mobile phone"
$wsContacts.Range("H3").Value = "Fax"
$wsContacts.Range("I3").Value = "Notes"

# ---------------------------------------------------------------------
# 4. Refresh the header-row fill color (used by the B3:M3 / B3:I3 ranges
#    on every sheet) from the old dark green to the new dark teal.
# ---------------------------------------------------------------------
$headerColor = 5521920  # RGB(0, 66, 84) = #004254

$wsSuppliers.Range("B3:M3").Interior.Color = $headerColor
$wsFactories.Range("B3:M3").Interior.Color = $headerColor
$wsContacts.Range("B3:I3").Interior.Color = $headerColor
